$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 102.8
$ws.Range("I15").Value = 102.8
$ws.Range("K15").Value = 308.4
$ws.Range("M15").Value = -139.4

$ws.Range("H17").Value = 2278813.2
$ws.Range("J17").Value = 2278813.2
$ws.Range("L17").Value = 6836439.600000001
$ws.Range("N17").Value = -6836775.600000001

$ws.Range("H64").Value = 3295.814
$ws.Range("I64").Value = 2966.2856
$ws.Range("J64").Value = 4737.5
$ws.Range("K64").Value = 2966.2856
$ws.Range("L64").Value = 4737.5
$ws.Range("M64").Value = -2718.2856
$ws.Range("N64").Value = -5233.5

$ws.Range("H67").Value = 3295.814
$ws.Range("I67").Value = 2966.2856
$ws.Range("J67").Value = 4737.5
$ws.Range("K67").Value = 2966.2856
$ws.Range("L67").Value = 4737.5
$ws.Range("M67").Value = -2108.2856
$ws.Range("N67").Value = -6453.5

$ws.Range("H69").Value = 3485.9524
$ws.Range("J69").Value = 3915.9092
$ws.Range("L69").Value = 11747.7276
$ws.Range("N69").Value = -13495.7276

$ws.Range("H72").Value = 3485.9524
$ws.Range("J72").Value = 3915.9092
$ws.Range("L72").Value = 35243.1828
$ws.Range("N72").Value = -43979.1828

$ws.Range("H74").Value = 4500
$ws.Range("I74").Value = 5000
$ws.Range("J74").Value = 4400
$ws.Range("L74").Value = 4400
$ws.Range("M74").Value = -4064
$ws.Range("N74").Value = -6272

$ws.Range("H77").Value = 4500
$ws.Range("I77").Value = 5000
$ws.Range("J77").Value = 4400
$ws.Range("K77").Value = 25000
$ws.Range("L77").Value = 22000
$ws.Range("M77").Value = -20320
$ws.Range("N77").Value = -31360

$ws.Range("H112").Value = 1285.6522
$ws.Range("I112").Value = 900
$ws.Range("J112").Value = 1322.381
$ws.Range("K112").Value = 2700
$ws.Range("L112").Value = 3967.143
$ws.Range("M112").Value = -1592
$ws.Range("N112").Value = -6183.143

$ws.Range("H137").Value = 1266.7715
$ws.Range("I137").Value = 1276.3914
$ws.Range("J137").Value = 1248.3334
$ws.Range("K137").Value = 3829.1742
$ws.Range("L137").Value = 3745.0002
$ws.Range("M137").Value = -1279.1742
$ws.Range("N137").Value = -8845.0002

$ws.Range("H138").Value = 6505.9155
$ws.Range("I138").Value = 3114.7368
$ws.Range("J138").Value = 7512.672
$ws.Range("K138").Value = 9344.2104
$ws.Range("L138").Value = 22538.016
$ws.Range("M138").Value = -4204.2104
$ws.Range("N138").Value = -32818.016

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15240.322
$ws.Range("I32").Value = 12825
$ws.Range("J32").Value = 27800
$ws.Range("K32").Value = 12825
$ws.Range("L32").Value = 27800
$ws.Range("M32").Value = -12538
$ws.Range("N32").Value = -28374

$ws.Range("H45").Value = 1972.32
$ws.Range("I45").Value = 1385.7368
$ws.Range("K45").Value = 1385.7368
$ws.Range("M45").Value = -1008.7368

$ws.Range("H122").Value = 2323.5386
$ws.Range("I122").Value = 2029.2
$ws.Range("J122").Value = 3304.6667
$ws.Range("K122").Value = 6087.6
$ws.Range("L122").Value = 9914.000100000001
$ws.Range("M122").Value = -3637.6
$ws.Range("N122").Value = -14814.0001

$ws.Range("H132").Value = 1701.8518
$ws.Range("I132").Value = 1340.1777
$ws.Range("J132").Value = 3510.2222
$ws.Range("K132").Value = 4020.5331
$ws.Range("L132").Value = 10530.6666
$ws.Range("M132").Value = -1490.5331
$ws.Range("N132").Value = -15590.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 911.3077
$ws.Range("I7").Value = 1282.1111
$ws.Range("K7").Value = 1282.1111
$ws.Range("M7").Value = -1169.1111

$ws.Range("H31").Value = 4232.7144
$ws.Range("I31").Value = 1776
$ws.Range("K31").Value = 1776
$ws.Range("M31").Value = -1481

$ws.Range("H34").Value = 4232.7144
$ws.Range("I34").Value = 1776
$ws.Range("K34").Value = 1776
$ws.Range("M34").Value = -1574

$ws.Range("H58").Value = 2501.5417
$ws.Range("I58").Value = 1667.1111
$ws.Range("J58").Value = 5004.8335
$ws.Range("K58").Value = 1667.1111
$ws.Range("L58").Value = 5004.8335
$ws.Range("M58").Value = -1464.1111
$ws.Range("N58").Value = -5410.8335

$ws.Range("H105").Value = 2799.8572
$ws.Range("I105").Value = 4150
$ws.Range("J105").Value = 999.6667
$ws.Range("K105").Value = 4150
$ws.Range("L105").Value = 999.6667
$ws.Range("M105").Value = -2403
$ws.Range("N105").Value = -4493.6667

$ws.Range("H136").Value = 2501.5417
$ws.Range("I136").Value = 1667.1111
$ws.Range("J136").Value = 5004.8335
$ws.Range("K136").Value = 5001.3333
$ws.Range("L136").Value = 15014.5005
$ws.Range("M136").Value = -2451.3333
$ws.Range("N136").Value = -20114.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 332706.28
$ws.Range("I107").Value = 956.8333
$ws.Range("J107").Value = 592336.3
$ws.Range("K107").Value = 2870.4999
$ws.Range("L107").Value = 1777008.9
$ws.Range("M107").Value = -950.4998999999998
$ws.Range("N107").Value = -1780848.9

$ws.Range("H113").Value = 1159.0526
$ws.Range("I113").Value = 1415.75
$ws.Range("J113").Value = 719
$ws.Range("K113").Value = 4247.25
$ws.Range("L113").Value = 2157
$ws.Range("M113").Value = -2077.25
$ws.Range("N113").Value = -6497

$ws.Range("H122").Value = 990.5484
$ws.Range("I122").Value = 591.8461
$ws.Range("J122").Value = 3063.8
$ws.Range("K122").Value = 5326.6149
$ws.Range("L122").Value = 27574.2
$ws.Range("M122").Value = -2876.6149
$ws.Range("N122").Value = -32474.2

$ws.Range("H132").Value = 1468.6957
$ws.Range("I132").Value = 411.81818
$ws.Range("J132").Value = 2437.5
$ws.Range("K132").Value = 3706.36362
$ws.Range("L132").Value = 21937.5
$ws.Range("M132").Value = -1176.36362
$ws.Range("N132").Value = -26997.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 936.9545000000001
$ws.Range("I2").Value = 735.9286
$ws.Range("J2").Value = 1288.75
$ws.Range("K2").Value = 735.9286
$ws.Range("L2").Value = 1288.75
$ws.Range("M2").Value = -622.9286
$ws.Range("N2").Value = -1514.75

$ws.Range("H80").Value = 2280.7144
$ws.Range("I80").Value = 1652.5
$ws.Range("J80").Value = 2532
$ws.Range("K80").Value = 1652.5
$ws.Range("L80").Value = 2532
$ws.Range("M80").Value = -654.5
$ws.Range("N80").Value = -4528

$ws.Range("H83").Value = 2280.7144
$ws.Range("I83").Value = 1652.5
$ws.Range("J83").Value = 2532
$ws.Range("K83").Value = 8262.5
$ws.Range("L83").Value = 12660
$ws.Range("M83").Value = -3270.5
$ws.Range("N83").Value = -22644

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1361.75
$ws.Range("I22").Value = 998.6667
$ws.Range("J22").Value = 1828.5714
$ws.Range("K22").Value = 998.6667
$ws.Range("L22").Value = 1828.5714
$ws.Range("M22").Value = -703.6667
$ws.Range("N22").Value = -2418.5714

$ws.Range("H27").Value = 1361.75
$ws.Range("I27").Value = 998.6667
$ws.Range("J27").Value = 1828.5714
$ws.Range("K27").Value = 998.6667
$ws.Range("L27").Value = 1828.5714
$ws.Range("M27").Value = -891.6667
$ws.Range("N27").Value = -2042.5714

$ws.Range("H133").Value = 40473
$ws.Range("J133").Value = 40473
$ws.Range("L133").Value = 40473
$ws.Range("N133").Value = -45533

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 10000
$ws.Range("J75").Value = 10000
$ws.Range("L75").Value = 10000
$ws.Range("N75").Value = -11872

$ws.Range("H78").Value = 10000
$ws.Range("J78").Value = 10000
$ws.Range("L78").Value = 30000
$ws.Range("N78").Value = -39360

$ws.Range("H81").Value = 2333.389
$ws.Range("I81").Value = 2000.125
$ws.Range("J81").Value = 2600
$ws.Range("K81").Value = 4000.25
$ws.Range("L81").Value = 5200
$ws.Range("M81").Value = -2939.25
$ws.Range("N81").Value = -7322

$ws.Range("H84").Value = 2333.389
$ws.Range("I84").Value = 2000.125
$ws.Range("J84").Value = 2600
$ws.Range("K84").Value = 20001.25
$ws.Range("L84").Value = 26000
$ws.Range("M84").Value = -14697.25
$ws.Range("N84").Value = -36608
